{"js": "// Update the date line and the five rows of division-practice answers.\n\n// 1) Update the date heading.\nconst dateResults = context.document.body.search(\"2024-10-09 Wednesday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2024-10-10 Thursday\", Word.InsertLocation.replace);\n}\n\n// 2) Update the table of answers. The table has 20 rows, but only every\n// 4th row (0, 4, 8, 12, 16) holds the 5 answer cells; the rows in between\n// are blank spacer rows. Target each cell by (row, column) so duplicate\n// text values elsewhere in the table aren't accidentally overwritten.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = {\n  0: [\"95\u00f73=31, 2\", \"93\u00f73=31, 0\", \"62\u00f75=12, 2\", \"56\u00f79=6, 2\", \"24\u00f79=2, 6\"],\n  4: [\"53\u00f79=5, 8\", \"66\u00f74=16, 2\", \"52\u00f73=17, 1\", \"32\u00f77=4, 4\", \"81\u00f72=40, 1\"],\n  8: [\"71\u00f73=23, 2\", \"24\u00f76=4, 0\", \"26\u00f76=4, 2\", \"34\u00f74=8, 2\", \"24\u00f79=2, 6\"],\n  12: [\"28\u00f77=4, 0\", \"86\u00f78=10, 6\", \"22\u00f74=5, 2\", \"39\u00f78=4, 7\", \"52\u00f78=6, 4\"],\n  16: [\"90\u00f72=45, 0\", \"95\u00f74=23, 3\", \"35\u00f76=5, 5\", \"92\u00f72=46, 0\", \"48\u00f72=24, 0\"],\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const row = parseInt(rowIndex, 10);\n  const rowValues = newValues[rowIndex];\n  for (let col = 0; col < rowValues.length; col++) {\n    table.getCell(row, col).value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the five rows of division-practice answers.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date heading.\n$find = $d.Content.Find\n$find.Text = \"2024-10-09 Wednesday\"\n$find.Replacement.Text = \"2024-10-10 Thursday\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, \"2024-10-10 Thursday\", 2)\n\n# 2) Update the table of answers. The table has 20 rows, but only every\n# 4th row (Word rows 1, 5, 9, 13, 17) holds the 5 answer cells; the rows\n# in between are blank spacer rows. Target each cell by (row, column) so\n# duplicate text values elsewhere in the table aren't accidentally\n# overwritten.\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"95\u00f73=31, 2\", \"93\u00f73=31, 0\", \"62\u00f75=12, 2\", \"56\u00f79=6, 2\", \"24\u00f79=2, 6\")\n    5  = @(\"53\u00f79=5, 8\", \"66\u00f74=16, 2\", \"52\u00f73=17, 1\", \"32\u00f77=4, 4\", \"81\u00f72=40, 1\")\n    9  = @(\"71\u00f73=23, 2\", \"24\u00f76=4, 0\", \"26\u00f76=4, 2\", \"34\u00f74=8, 2\", \"24\u00f79=2, 6\")\n    13 = @(\"28\u00f77=4, 0\", \"86\u00f78=10, 6\", \"22\u00f74=5, 2\", \"39\u00f78=4, 7\", \"52\u00f78=6, 4\")\n    17 = @(\"90\u00f72=45, 0\", \"95\u00f74=23, 3\", \"35\u00f76=5, 5\", \"92\u00f72=46, 0\", \"48\u00f72=24, 0\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowValues = $newValues[$rowIndex]\n    for ($c = 1; $c -le $rowValues.Length; $c++) {\n        $t.Cell($rowIndex, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
